$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates (shared strings): issue number and the reporting week
# date range.
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  10"
$ws.Range("C9").Value = "Report Covering the Week  3/6/2023  Through  3/12/2023"

# ---------------------------------------------------------------------------
# Cells that flip from a numeric value to the literal text "0" or "***.*"
# (used whenever one side of a % change is zero/undefined). A14's sibling
# cell D14 already carries the exact target style (s=14, General format) and
# is never touched by this edit, so it is used as the formats-only paste
# source to keep the style index identical to the target workbook.
# L15 already carries the "***.*" style (s=14) and is likewise untouched.
# ---------------------------------------------------------------------------
function Set-TextCell($ref, $text, $formatSource) {
    $ws.Range($ref).Value = "'" + $text
    $ws.Range($formatSource).Copy() | Out-Null
    $ws.Range($ref).PasteSpecial(-4122) | Out-Null
}

Set-TextCell "C14" "0" "D14"
Set-TextCell "D15" "0" "C15"
Set-TextCell "E15" "***.*" "L15"
Set-TextCell "C17" "0" "D14"
Set-TextCell "C22" "0" "D14"
Set-TextCell "D26" "0" "C26"
Set-TextCell "E26" "***.*" "L15"
Set-TextCell "D27" "0" "C26"
Set-TextCell "E27" "***.*" "L15"
Set-TextCell "C28" "0" "D28"
Set-TextCell "C29" "0" "D29"

$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Plain numeric value updates.
# ---------------------------------------------------------------------------
$ws.Range("N15").Value = -60

$ws.Range("D16").Value = 7
$ws.Range("E16").Value = -57.142857142857
$ws.Range("F16").Value = 14
$ws.Range("G16").Value = 23
$ws.Range("H16").Value = -39.130434782608
$ws.Range("I16").Value = 35
$ws.Range("J16").Value = 44
$ws.Range("K16").Value = -20.454545454545
$ws.Range("L16").Value = 9.375
$ws.Range("M16").Value = 40
$ws.Range("N16").Value = -86.988847583643

$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -100
$ws.Range("G17").Value = 12
$ws.Range("H17").Value = -25
$ws.Range("I17").Value = 32
$ws.Range("J17").Value = 31
$ws.Range("K17").Value = 3.225806451612
$ws.Range("L17").Value = 68.421052631578
$ws.Range("M17").Value = 88.235294117647
$ws.Range("N17").Value = -30.434782608695

$ws.Range("C18").Value = 8
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 300
$ws.Range("F18").Value = 27
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = 350
$ws.Range("I18").Value = 55
$ws.Range("J18").Value = 47
$ws.Range("K18").Value = 17.021276595744
$ws.Range("L18").Value = 17.021276595744
$ws.Range("M18").Value = -23.611111111111
$ws.Range("N18").Value = -89.382239382239

$ws.Range("C19").Value = 37
$ws.Range("D19").Value = 28
$ws.Range("E19").Value = 32.142857142857
$ws.Range("F19").Value = 130
$ws.Range("G19").Value = 123
$ws.Range("H19").Value = 5.691056910569
$ws.Range("I19").Value = 296
$ws.Range("J19").Value = 291
$ws.Range("K19").Value = 1.718213058419
$ws.Range("L19").Value = 58.288770053475
$ws.Range("M19").Value = 23.849372384937
$ws.Range("N19").Value = -59.618008185538

$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 16.666666666666
$ws.Range("I20").Value = 19
$ws.Range("J20").Value = 22
$ws.Range("K20").Value = -13.636363636363
$ws.Range("L20").Value = 18.75
$ws.Range("M20").Value = 171.428571428571
$ws.Range("N20").Value = -97.342657342657

$ws.Range("C21").Value = 49
$ws.Range("D21").Value = 43
$ws.Range("E21").Value = 13.953488372093
$ws.Range("F21").Value = 188
$ws.Range("G21").Value = 171
$ws.Range("H21").Value = 9.941520467836
$ws.Range("I21").Value = 440
$ws.Range("J21").Value = 437
$ws.Range("K21").Value = 0.686498855835
$ws.Range("L21").Value = 46.179401993355
$ws.Range("M21").Value = 21.212121212121
$ws.Range("N21").Value = -80.760822037603

$ws.Range("E22").Value = -100
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = -25
$ws.Range("J22").Value = 11
$ws.Range("K22").Value = -36.363636363636

$ws.Range("L23").Value = -66.666666666666
$ws.Range("M23").Value = -71.428571428571

$ws.Range("C24").Value = 55
$ws.Range("D24").Value = 94
$ws.Range("E24").Value = -41.489361702127
$ws.Range("F24").Value = 246
$ws.Range("G24").Value = 296
$ws.Range("H24").Value = -16.891891891891
$ws.Range("I24").Value = 637
$ws.Range("J24").Value = 623
$ws.Range("K24").Value = 2.247191011235
$ws.Range("L24").Value = 43.468468468468
$ws.Range("M24").Value = 113.758389261745

$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = -30
$ws.Range("F25").Value = 29
$ws.Range("H25").Value = -14.705882352941
$ws.Range("I25").Value = 77
$ws.Range("J25").Value = 84
$ws.Range("K25").Value = -8.333333333333
$ws.Range("L25").Value = 45.283018867924
$ws.Range("M25").Value = 16.666666666666

$ws.Range("F27").Value = 7
$ws.Range("G27").Value = 7
$ws.Range("I27").Value = 21
$ws.Range("K27").Value = 5
$ws.Range("L27").Value = 133.333333333333

$ws.Range("G30").Value = 9
$ws.Range("J30").Value = 11
$ws.Range("K30").Value = -90.909090909090
